# Updated cryptos list with latest price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.225.88"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "1.572.58"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "1.796.24"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "1.573.19"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.191.83"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0703"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "216.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.76%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "1.454.75"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "1.707.37"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").Value = "  +3.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0523"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.75%  "
